$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.964.22'
$ws.Range('E2').Value = '  +2.23%  '
$ws.Range('D3').Value = '2.052.15'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '229.23'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.617'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '58.68'
$ws.Range('E7').Value = '  +6.54%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +2.07%  '
$ws.Range('E10').Value = '  +2.97%  '
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').Value = '2.353.49'
$ws.Range('E12').Value = '  +1.72%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.66'
$ws.Range('E13').Value = '  +3.27%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.85'
$ws.Range('E14').Value = '  +3.15%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.753'
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('E16').Value = '  +1.92%  '
$ws.Range('D17').Value = '2.039.02'
$ws.Range('E17').Value = '  +0.86%  '
$ws.Range('D18').Value = '37.897.50'
$ws.Range('E18').Value = '  +2.21%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.29'
$ws.Range('E19').Value = '  -3.33%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.62'
$ws.Range('D21').Value = '0.0₃0836'
$ws.Range('E21').Value = '  +2.53%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '224.55'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  +2.21%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.30'
$ws.Range('E26').Value = '  +0.96%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '166.44'
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.134'
$ws.Range('E28').Value = '  +4.64%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.04'
$ws.Range('E29').Value = '  +1.71%  '
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.52'
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.58'
$ws.Range('E33').Value = '  +2.39%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.06'
$ws.Range('E34').Value = '  +10.28%  '
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('E36').Value = '  -0.65%  '
$ws.Range('E37').Value = '  +10.03%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.28'
$ws.Range('E38').Value = '  +5.48%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '1.492.48'
$ws.Range('E40').Value = '  +1.58%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0217'
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '97.23'
$ws.Range('E42').Value = '  +1.69%  '
$ws.Range('E43').Value = '  +2.88%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.52'
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0919'
$ws.Range('E45').Value = '  +1.10%  '
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.15'
$ws.Range('E47').Value = '  +13.77%  '
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('E50').Value = '  -2.46%  '
$ws.Range('D51').Value = '2.243.07'
$ws.Range('E51').Value = '  +1.55%  '
